$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.984.35'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.760.19'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3401'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.40'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.119'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.25'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.150'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.77%  '
$ws.Range("D15").Value = '1.758.12'
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.098'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001058'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06629'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.35'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9984'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.93'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.226'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("D23").Value = '28.001.76'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.63'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.389'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.80'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.94'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.311'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").Value = '1.958.23'
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.278'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -11.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.05'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.076'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.817'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08710'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.07'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.59%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02288'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06168'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.19%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.133'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6494'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2109'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.500'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.205'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9986'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.863'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.77'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.831'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("E47").Value = '  -4.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.75'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.002'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.86%  '
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06996'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.09%  '
